$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# --- Simple value corrections (existing rows) ---
$ws.Range("F8").Value = 15000
$ws.Range("F40").Value = 4000

# --- Make room for a new sale record between the current row 63 (JAI SRI ...)
#     and row 64 (SRI MANJUNATHA ...): shift old rows 64-70 down to 65-71 ---
for ($r = 70; $r -ge 64; $r--) {
    $src = $ws.Range("A" + $r + ":H" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":H" + ($r + 1))
    $src.Copy()
    $dst.PasteSpecial(-4104)
}
$excel.CutCopyMode = 0
# Match the formatting of the row directly above the new blank row
$ws.Range("A63:H63").Copy()
$ws.Range("A64:H64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new record at row 64
$ws.Range("A64").Value = "JAI SRI ELECTRICALS & HARDWARE"
$ws.Range("B64").Value = "Dasanpura"
$ws.Range("C64").Value = 126
$ws.Range("D64").Value = 45283
$ws.Range("E64").Value = 23800
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0

# --- Make room for two more new sale records between the now-shifted row 68
#     (NAVARTHNA ELECTRICALS) and row 69 (Chandre gowda ...): shift rows
#     69-71 down to 71-73 ---
for ($r = 71; $r -ge 69; $r--) {
    $src = $ws.Range("A" + $r + ":H" + $r)
    $dst = $ws.Range("A" + ($r + 2) + ":H" + ($r + 2))
    $src.Copy()
    $dst.PasteSpecial(-4104)
}
$excel.CutCopyMode = 0
# Match the formatting of the row directly above the new blank rows
$ws.Range("A68:H68").Copy()
$ws.Range("A69:H70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new record at row 69
$ws.Range("A69").Value = "R B LIGHTING"
$ws.Range("B69").Value = "SHIVANAGAR "
$ws.Range("C69").Value = 127
$ws.Range("D69").Value = 45285
$ws.Range("E69").Value = 19100
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0

# Fill in the new record at row 70
$ws.Range("A70").Value = "BRIGHT LIGHTS"
$ws.Range("B70").Value = "MAGADI MAIN ROAD "
$ws.Range("C70").Value = 128
$ws.Range("D70").Value = 45286
$ws.Range("E70").Value = 18600
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0

# --- View/selection tidy-up to match what was saved with the workbook ---
$ws.Activate()
$ws.Range("F8").Select()
